$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 305.6
$ws.Range("I11").Value = 305.6
$ws.Range("K11").Value = 305.6
$ws.Range("M11").Value = -165.6
$ws.Range("H12").Value = 205.4
$ws.Range("J12").Value = 225
$ws.Range("L12").Value = 225
$ws.Range("N12").Value = -565
$ws.Range("H40").Value = 2539.2
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2325
$ws.Range("H51").Value = 2500
$ws.Range("I51").Value = 2500
$ws.Range("K51").Value = 2500
$ws.Range("M51").Value = -2016
$ws.Range("H55").Value = 428.33334
$ws.Range("I55").Value = 341
$ws.Range("K55").Value = 341
$ws.Range("M55").Value = -127
$ws.Range("H62").Value = 3970.25
$ws.Range("I62").Value = 3664.3333
$ws.Range("J62").Value = 4888
$ws.Range("K62").Value = 3664.3333
$ws.Range("L62").Value = 4888
$ws.Range("M62").Value = -3040.3333
$ws.Range("N62").Value = -6136
$ws.Range("H65").Value = 3970.25
$ws.Range("I65").Value = 3664.3333
$ws.Range("J65").Value = 4888
$ws.Range("K65").Value = 18321.6665
$ws.Range("L65").Value = 24440
$ws.Range("M65").Value = -15201.6665
$ws.Range("N65").Value = -30680
$ws.Range("H69").Value = 6000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 6000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 18000
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 6000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 6000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 54000
$ws.Range("N72").Value = -62736
$ws.Range("H80").Value = 2171.2856
$ws.Range("I80").Value = 675
$ws.Range("J80").Value = 4166.3335
$ws.Range("K80").Value = 2025
$ws.Range("L80").Value = 12499.0005
$ws.Range("M80").Value = -1027
$ws.Range("N80").Value = -14495.0005
$ws.Range("H83").Value = 2171.2856
$ws.Range("I83").Value = 675
$ws.Range("J83").Value = 4166.3335
$ws.Range("K83").Value = 6075
$ws.Range("L83").Value = 37497.0015
$ws.Range("M83").Value = -1083
$ws.Range("N83").Value = -47481.0015
$ws.Range("H98").Value = 988.9091
$ws.Range("I98").Value = 1019.7778
$ws.Range("K98").Value = 1019.7778
$ws.Range("M98").Value = 478.2222
$ws.Range("H100").Value = 2695.3572
$ws.Range("I100").Value = 2936.25
$ws.Range("J100").Value = 1250
$ws.Range("K100").Value = 2936.25
$ws.Range("L100").Value = 1250
$ws.Range("M100").Value = -2395.25
$ws.Range("N100").Value = -2332
$ws.Range("H107").Value = 3003.4
$ws.Range("I107").Value = 3003.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3003.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1083.4
$ws.Range("H121").Value = 2999.25
$ws.Range("J121").Value = 2999.25
$ws.Range("L121").Value = 8997.75
$ws.Range("N121").Value = -12491.75
$ws.Range("H122").Value = 988.9091
$ws.Range("I122").Value = 1019.7778
$ws.Range("K122").Value = 3059.3334
$ws.Range("M122").Value = -609.3334
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4572.1333
$ws.Range("I2").Value = 1673.8889
$ws.Range("K2").Value = 1673.8889
$ws.Range("M2").Value = -1560.8889
$ws.Range("H32").Value = 3765.8667
$ws.Range("I32").Value = 2224.4583
$ws.Range("K32").Value = 2224.4583
$ws.Range("M32").Value = -1937.4583
$ws.Range("H97").Value = 1093.4
$ws.Range("I97").Value = 611.13336
$ws.Range("K97").Value = 611.13336
$ws.Range("M97").Value = -115.13336
$ws.Range("H116").Value = 4572.1333
$ws.Range("I116").Value = 1673.8889
$ws.Range("K116").Value = 1673.8889
$ws.Range("M116").Value = 620.1111000000001
$ws.Range("H122").Value = 2982.4666
$ws.Range("I122").Value = 2666
$ws.Range("K122").Value = 7998
$ws.Range("M122").Value = -5548

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4572.1333
$ws.Range("I3").Value = 1673.8889
$ws.Range("K3").Value = 1673.8889
$ws.Range("M3").Value = -1559.8889
$ws.Range("H10").Value = 999.5
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 999.5
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 999.5
$ws.Range("N10").Value = -1279.5
$ws.Range("H22").Value = 290.83334
$ws.Range("I22").Value = 211.25
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 211.25
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -38.25
$ws.Range("N22").Value = -796
$ws.Range("H105").Value = 3680.4375
$ws.Range("I105").Value = 3337.5386
$ws.Range("K105").Value = 3337.5386
$ws.Range("M105").Value = -1590.5386
$ws.Range("H134").Value = 1683.7273
$ws.Range("I134").Value = 1552.1
$ws.Range("K134").Value = 4656.299999999999
$ws.Range("M134").Value = -2121.299999999999
$ws.Range("M10").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("H6").Value = 171.5
$ws.Range("I6").Value = 101
$ws.Range("K6").Value = 101
$ws.Range("M6").Value = 12
$ws.Range("H10").Value = 427.7143
$ws.Range("I10").Value = 415.66666
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 415.66666
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = -276.66666
$ws.Range("N10").Value = -778
$ws.Range("H31").Value = 2122.625
$ws.Range("I31").Value = 1911.3334
$ws.Range("K31").Value = 1911.3334
$ws.Range("M31").Value = -1616.3334
$ws.Range("H34").Value = 2122.625
$ws.Range("I34").Value = 1911.3334
$ws.Range("K34").Value = 1911.3334
$ws.Range("M34").Value = -1709.3334
$ws.Range("H62").Value = 3001
$ws.Range("I62").Value = 3001
$ws.Range("K62").Value = 3001
$ws.Range("M62").Value = -2377
$ws.Range("H65").Value = 3001
$ws.Range("I65").Value = 3001
$ws.Range("K65").Value = 15005
$ws.Range("M65").Value = -11885
$ws.Range("H105").Value = 1486.5834
$ws.Range("I105").Value = 1613.5714
$ws.Range("J105").Value = 1308.8
$ws.Range("K105").Value = 1613.5714
$ws.Range("L105").Value = 1308.8
$ws.Range("M105").Value = 133.4286
$ws.Range("N105").Value = -4802.8
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 999
$ws.Range("J75").Value = 999
$ws.Range("L75").Value = 2997
$ws.Range("N75").Value = -4993
$ws.Range("H78").Value = 999
$ws.Range("J78").Value = 999
$ws.Range("L78").Value = 8991
$ws.Range("N78").Value = -18975
$ws.Range("H131").Value = 1125.3846
$ws.Range("J131").Value = 1127.5
$ws.Range("L131").Value = 3382.5
$ws.Range("N131").Value = -13462.5
$ws.Range("H132").Value = 2399
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 6699.4
$ws.Range("J19").Value = 6699.4
$ws.Range("L19").Value = 6699.4
$ws.Range("N19").Value = -7275.4
$ws.Range("H102").Value = 1428.8572
$ws.Range("J102").Value = 1665.3334
$ws.Range("L102").Value = 1665.3334
$ws.Range("N102").Value = -4909.3334
$ws.Range("H122").Value = 1628.8572
$ws.Range("J122").Value = 988
$ws.Range("L122").Value = 2964
$ws.Range("N122").Value = -7864

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1990
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1990
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1990
$ws.Range("N7").Value = -2214
$ws.Range("H12").Value = 1350
$ws.Range("I12").Value = 750
$ws.Range("K12").Value = 750
$ws.Range("M12").Value = -580
$ws.Range("H22").Value = 1096
$ws.Range("I22").Value = 999.5
$ws.Range("K22").Value = 999.5
$ws.Range("M22").Value = -704.5
$ws.Range("H27").Value = 1096
$ws.Range("I27").Value = 999.5
$ws.Range("K27").Value = 999.5
$ws.Range("M27").Value = -892.5
$ws.Range("H46").Value = 1603.6786
$ws.Range("I46").Value = 1118.6875
$ws.Range("K46").Value = 1118.6875
$ws.Range("M46").Value = -930.6875
$ws.Range("H126").Value = 1990
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1990
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 5970
$ws.Range("N126").Value = -10910
$ws.Range("H136").Value = 4004
$ws.Range("I136").Value = 4004
$ws.Range("K136").Value = 12012
$ws.Range("M136").Value = -9462
$ws.Range("M7").ClearContents()
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 4500
$ws.Range("J8").Value = 4500
$ws.Range("L8").Value = 4500
$ws.Range("N8").Value = -4780
$ws.Range("H13").Value = 5000
$ws.Range("J13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("N13").Value = -5280
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H132").Value = 3528.2856
$ws.Range("I132").Value = 2897.6667
$ws.Range("K132").Value = 8693.000100000001
$ws.Range("M132").Value = -6163.000100000001
$ws.Range("H136").Value = 1394.6428
$ws.Range("I136").Value = 1075.1052
$ws.Range("J136").Value = 2069.2222
$ws.Range("K136").Value = 3225.3156
$ws.Range("L136").Value = 6207.6666
$ws.Range("M136").Value = -675.3155999999999
$ws.Range("N136").Value = -11307.6666
$ws.Range("N49").ClearContents()
